$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# New trainees added to the roster: row 17 (Chandrashekar S / user15)
# and row 18 (Indiver Jamwal / user16). Row 19 stays the blank filler
# row that was already there.
# ------------------------------------------------------------------

# --- Row 17 ---------------------------------------------------------
$ws.Range("B17").Value = 15

# C17 needs a "boxed" look (thin border on left+right only) that isn't
# used anywhere else yet, so build it by hand on the cell itself.
$ws.Range("C17").Value = "Chandrashekar S"
$c17 = $ws.Range("C17")
$c17.Font.Name = $ws.Range("D3").Font.Name
$c17.Borders.Item(7).LineStyle = 1
$c17.Borders.Item(7).Weight = 2
$c17.Borders.Item(10).LineStyle = 1
$c17.Borders.Item(10).Weight = 2
$c17.VerticalAlignment = -4108

# D17 reuses the same format as the existing "Server" column cells.
$ws.Range("D3").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "http://34.238.115.168:8000/en-US/app/launcher/home"

# E17 is another new thin left/right boxed style (plain font this time).
$ws.Range("E17").Value = "user15"
$e17 = $ws.Range("E17")
$e17.Borders.Item(7).LineStyle = 1
$e17.Borders.Item(7).Weight = 2
$e17.Borders.Item(10).LineStyle = 1
$e17.Borders.Item(10).Weight = 2

# F17 carries a hyperlink back to the Splunk server, but keeps the
# plain boxed look used by the rest of column F.
$ws.Range("F17").Value = "user@123"
$ws.Hyperlinks.Add($ws.Range("F17"), "http://34.238.115.168:8000/en-US/app/launcher/home") | Out-Null
$ws.Range("F16").Copy()
$ws.Range("F17").PasteSpecial(-4122)

# --- Row 18 ----------------------------------------------------------
$ws.Range("B18").Value = 16

# C18 was already a blank filler cell sharing the plain "no border"
# style used elsewhere, so just fill in the value.
$ws.Range("C18").Value = "Indiver Jamwal"

# D18 reuses the alternating grey-row format used further up the table.
$ws.Range("D4").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "http://174.129.191.6:8000/en-US/app/launcher/home"

# E18 uses the same new boxed style as E17.
$ws.Range("E17").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = "user16"

# F18 also gets a hyperlink, formatted like the rest of column F.
$ws.Range("F18").Value = "user@123"
$ws.Hyperlinks.Add($ws.Range("F18"), "http://174.129.191.6:8000/en-US/app/launcher/home") | Out-Null
$ws.Range("F16").Copy()
$ws.Range("F18").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Move the active selection to where the user last left off editing.
$ws.Range("D14").Select()
